$d = $word.ActiveDocument

# Find the "Docente(s) Responsável(eis)" heading paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Docente\(s\) Responsável\(eis\)") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Docente(s) Responsável(eis)' heading paragraph"
}

# Insert a new paragraph right after the heading, then turn it into a
# bullet-list item naming the instructor.
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$newPara = $target.Next()
$newPara.Style = "ListBullet"
$newPara.Range.Text = "6712818 - Mauricio Lamano Ferreira"

Write-Output "Inserted '$($newPara.Range.Text)' after paragraph $i"
